# Add a new "ProductDesc" worksheet after the last existing sheet (Shipping)
# and populate it with a small product description table, matching the
# "Added a test for pu method on Product description" commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ProductDesc"

# Populate the data in the same column order used when the shared strings
# table was originally built (description, name, slug, then id) so the
# resulting shared-string ordering matches the source data.
$ws.Range("D1").Value = "description"
$ws.Range("D2").Value = "Curabitur non nulla sit amet nisl tempus convallis quis ac lectus. Sed porttitor lectus nibh. Vestibulum ante ipsum primis in faucibus orci luctus et ultrices posuere cubilia Curae; Donec velit neque, auctor sit amet aliquam vel, ullamcorper sit amet ligula. Proin eget tortor risus. Cras ultricies ligula sed magna dictum porta. Quisque velit nisi, pretium ut lacinia in, elementum id enim. Vivamus suscipit tortor eget felis porttitor volutpat. Donec rutrum congue leo eget malesuada. Cras ultricies ligula sed magna dictum porta. Nulla quis lorem ut libero malesuada feugiat."

$ws.Range("B1").Value = "name"
$ws.Range("B2").Value = "V-Neck T-Shirt"

$ws.Range("C1").Value = "slug"
$ws.Range("C2").Value = "v-neck-t-shirt"

$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 11

# Make the selection on the new sheet match the authored workbook state.
$ws.Range("C2").Select()
